# Update the "想去人数" (want-to-go count) figures in column F for the
# "展览" and "全部类型" worksheets, as captured at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 265
$ws1.Range("F9").Value = 4958
$ws1.Range("F23").Value = 75
$ws1.Range("F27").Value = 38
$ws1.Range("F28").Value = 3227
$ws1.Range("F30").Value = 2551
$ws1.Range("F32").Value = 1596
$ws1.Range("F33").Value = 3729
$ws1.Range("F38").Value = 8
$ws1.Range("F39").Value = 947
$ws1.Range("F40").Value = 1209
$ws1.Range("F41").Value = 32
$ws1.Range("F43").Value = 589
$ws1.Range("F44").Value = 349
$ws1.Range("F47").Value = 3515

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 265
$ws4.Range("F10").Value = 4958
$ws4.Range("F26").Value = 75
$ws4.Range("F29").Value = 3227
$ws4.Range("F32").Value = 2551
$ws4.Range("F33").Value = 1596
$ws4.Range("F34").Value = 3729
$ws4.Range("F38").Value = 947
$ws4.Range("F40").Value = 1209
$ws4.Range("F41").Value = 32
$ws4.Range("F43").Value = 589
$ws4.Range("F49").Value = 3515
